{"js": "// Lab04 report tweak: after the short \"Find pins for USART3\" paragraph\n// (and before the bulleted \"Choose a set of RX/TX pins...\" item), add\n// three blank paragraphs using the same body style/indent so the list\n// has some breathing room for hand-written notes.\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\n// Locate the exact paragraph that precedes the new blank lines.\nlet anchor = null;\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  if (paragraphs.items[i].text === \"Find pins for USART3\") {\n    anchor = paragraphs.items[i];\n    break;\n  }\n}\nif (!anchor) {\n  throw new Error('Could not find the \"Find pins for USART3\" paragraph.');\n}\n\n// Build the three new paragraphs as raw OOXML so they come out as truly\n// empty paragraphs (pPr only, no run) instead of picking up a stray\n// empty run the way Paragraph.insertParagraph() would.\nconst blankPara =\n  '<w:p><w:pPr><w:pStyle w:val=\"CambraiMathLectureNormal\"/>' +\n  '<w:ind w:left=\"360\"/></w:pPr></w:p>';\n\nconst flatOpcXml =\n  '<?xml version=\"1.0\" encoding=\"UTF-8\" standalone=\"yes\"?>' +\n  '<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">' +\n  '<pkg:part pkg:name=\"/word/document.xml\" ' +\n  'pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\">' +\n  '<pkg:xmlData>' +\n  '<w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\">' +\n  \"<w:body>\" +\n  blankPara +\n  blankPara +\n  blankPara +\n  \"</w:body>\" +\n  \"</w:document>\" +\n  \"</pkg:xmlData></pkg:part></pkg:package>\";\n\nconst insertionPoint = anchor.getRange(Word.RangeLocation.after);\ninsertionPoint.insertOoxml(flatOpcXml, Word.InsertLocation.after);\n\nawait context.sync();\n", "ps1": "# Lab04 report tweak: after the short \"Find pins for USART3\" paragraph\n# (and before the bulleted \"Choose a set of RX/TX pins...\" item), add\n# three blank paragraphs using the same body style/indent so the list\n# has some breathing room for hand-written notes.\n\n$d = $word.ActiveDocument\n\n# Locate the exact paragraph that precedes the new blank lines.\n$findRange = $d.Content\n$found = $findRange.Find.Execute(\"Find pins for USART3\")\nif (-not $found) {\n    throw 'Could not find the \"Find pins for USART3\" paragraph.'\n}\n\n# Paragraphs(1) on the hit range expands it to the whole paragraph,\n# including its trailing paragraph mark, so .End sits right after it.\n$anchorParagraph = $findRange.Paragraphs(1)\n$insertionPoint = $d.Range($anchorParagraph.Range.End, $anchorParagraph.Range.End)\n\n# Insert as raw WordprocessingML so the new paragraphs end up truly\n# empty (pPr only, no run) instead of picking up a stray empty run.\n$blankParagraphXml = '<w:p xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\"><w:pPr><w:pStyle w:val=\"CambraiMathLectureNormal\"/><w:ind w:left=\"360\"/></w:pPr></w:p>'\n\n$insertionPoint.InsertXML($blankParagraphXml + $blankParagraphXml + $blankParagraphXml)\n"}
